$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.03677893029434
$ws.Range("D2").Value = 3.447146830702535
$ws.Range("E2").Value = 32.49987040633469
$ws.Range("F2").Value = 16.68346979744629
$ws.Range("G2").Value = 15.94996336549152
$ws.Range("H2").Value = 10.87453842930696
$ws.Range("L2").Value = 10.40915236689133
$ws.Range("M2").Value = 16.49099043153882
$ws.Range("O2").Value = 14.83260730414646
$ws.Range("B3").Value = 17.63493873695894
$ws.Range("D3").Value = 3.355693383895902
$ws.Range("E3").Value = 31.55155742348665
$ws.Range("F3").Value = 16.70905144835865
$ws.Range("G3").Value = 15.93902596904901
$ws.Range("H3").Value = 10.93857905381039
$ws.Range("L3").Value = 10.06436276062644
$ws.Range("M3").Value = 16.18472132630933
$ws.Range("O3").Value = 14.9261271388832
$ws.Range("B4").Value = 17.38404872809779
$ws.Range("D4").Value = 3.297904594803411
$ws.Range("E4").Value = 30.96035808651545
$ws.Range("F4").Value = 16.73339090180338
$ws.Range("G4").Value = 15.94599137398609
$ws.Range("H4").Value = 10.98090109865756
$ws.Range("L4").Value = 9.845504625723299
$ws.Range("M4").Value = 15.9934835715435
$ws.Range("O4").Value = 14.98996745018924
$ws.Range("B5").Value = 17.2808813888116
$ws.Range("D5").Value = 3.273968540141993
$ws.Range("E5").Value = 30.71759510861621
$ws.Range("F5").Value = 16.74546479901991
$ws.Range("G5").Value = 15.95224097036025
$ws.Range("H5").Value = 10.99889981839715
$ws.Range("L5").Value = 9.754616559096728
$ws.Range("M5").Value = 15.91482645286259
$ws.Range("O5").Value = 15.01758326395038
$ws.Range("B6").Value = 17.26369806016386
$ws.Range("D6").Value = 3.269971381843583
$ws.Range("E6").Value = 30.67718535381881
$ws.Range("F6").Value = 16.74759934961262
$ws.Range("G6").Value = 15.95348373025404
$ws.Range("H6").Value = 11.00193384885982
$ws.Range("L6").Value = 9.739424944229757
$ws.Range("M6").Value = 15.90172389111532
$ws.Range("O6").Value = 15.02226510414634
$ws.Range("B7").Value = 17.38266097331408
$ws.Range("D7").Value = 3.297583316349911
$ws.Range("E7").Value = 30.9570910181922
$ws.Range("F7").Value = 16.73354502945835
$ws.Range("G7").Value = 15.9460618910247
$ws.Range("H7").Value = 10.98114079351757
$ws.Range("L7").Value = 9.844285630734696
$ws.Range("M7").Value = 15.99242561303065
$ws.Range("O7").Value = 14.990333425115
$ws.Range("B8").Value = 17.89915664803664
$ws.Range("D8").Value = 3.415967509218214
$ws.Range("E8").Value = 32.17497822928038
$ws.Range("F8").Value = 16.69049024870109
$ws.Range("G8").Value = 15.94333769941106
$ws.Range("H8").Value = 10.89599533638305
$ws.Range("L8").Value = 10.29180978343397
$ws.Range("M8").Value = 16.38609467285305
$ws.Range("O8").Value = 14.86351255662421
$ws.Range("B9").Value = 18.8741123336738
$ws.Range("D9").Value = 3.634144013087377
$ws.Range("E9").Value = 34.47614994627365
$ws.Range("F9").Value = 16.67512001223808
$ws.Range("G9").Value = 16.04741812880989
$ws.Range("H9").Value = 10.75294563141905
$ws.Range("L9").Value = 11.10881737404024
$ws.Range("M9").Value = 17.12963029303975
$ws.Range("O9").Value = 14.66634956576953
$ws.Range("B10").Value = 19.56119870860157
$ws.Range("D10").Value = 3.784717520926388
$ws.Range("E10").Value = 36.09400822078517
$ws.Range("F10").Value = 16.70652904042688
$ws.Range("G10").Value = 16.19114508942749
$ws.Range("H10").Value = 10.66258112404143
$ws.Range("L10").Value = 11.66768891996687
$ws.Range("M10").Value = 17.65456277452841
$ws.Range("O10").Value = 14.55371774869877
$ws.Range("B11").Value = 19.86626626567164
$ws.Range("D11").Value = 3.850888050089113
$ws.Range("E11").Value = 36.81061766383552
$ws.Range("F11").Value = 16.73016267722024
$ws.Range("G11").Value = 16.27106345316184
$ws.Range("H11").Value = 10.62470468038988
$ws.Range("L11").Value = 11.91217820994386
$ws.Range("M11").Value = 17.88795625862116
$ws.Range("O11").Value = 14.50964328929621
$ws.Range("B12").Value = 19.98062221398386
$ws.Range("D12").Value = 3.875594632556466
$ws.Range("E12").Value = 37.07893439306032
$ws.Range("F12").Value = 16.74045765314679
$ws.Range("G12").Value = 16.30339765375837
$ws.Range("H12").Value = 10.61082932205434
$ws.Range("L12").Value = 12.00330467732381
$ws.Range("M12").Value = 17.97549958677848
$ws.Range("O12").Value = 14.49399668716661
$ws.Range("B13").Value = 19.95604684397648
$ws.Range("D13").Value = 3.870289470137818
$ws.Range("E13").Value = 37.02128703921907
$ws.Range("F13").Value = 16.73818061300317
$ws.Range("G13").Value = 16.29634220421456
$ws.Range("H13").Value = 10.6137967819967
$ws.Range("L13").Value = 11.98374453649751
$ws.Range("M13").Value = 17.95668375402569
$ws.Range("O13").Value = 14.49731984488934
$ws.Range("B14").Value = 19.87569826400158
$ws.Range("D14").Value = 3.852927792745704
$ws.Range("E14").Value = 36.83275473588773
$ws.Range("F14").Value = 16.73098268785881
$ws.Range("G14").Value = 16.27368224625068
$ws.Range("H14").Value = 10.6235537567823
$ws.Range("L14").Value = 11.91970471716685
$ws.Range("M14").Value = 17.89517562620865
$ws.Range("O14").Value = 14.50833504591817
$ws.Range("B15").Value = 19.8263279918799
$ws.Range("D15").Value = 3.842247129256678
$ws.Range("E15").Value = 36.7168689994589
$ws.Range("F15").Value = 16.7267489632572
$ws.Range("G15").Value = 16.26007132974697
$ws.Range("H15").Value = 10.62959117552082
$ws.Range("L15").Value = 11.88028730566428
$ws.Range("M15").Value = 17.85738927919841
$ws.Range("O15").Value = 14.51521845519053
$ws.Range("B16").Value = 19.54110329318442
$ws.Range("D16").Value = 3.780344886209456
$ws.Range("E16").Value = 36.04676334520144
$ws.Range("F16").Value = 16.70517278037736
$ws.Range("G16").Value = 16.18621327829852
$ws.Range("H16").Value = 10.66512171017147
$ws.Range("L16").Value = 11.65151001973068
$ws.Range("M16").Value = 17.63919597277744
$ws.Range("O16").Value = 14.55674340598287
$ws.Range("B17").Value = 19.36414275790876
$ws.Range("D17").Value = 3.741761999447674
$ws.Range("E17").Value = 35.63051991059118
$ws.Range("F17").Value = 16.69433211150348
$ws.Range("G17").Value = 16.14461576306599
$ws.Range("H17").Value = 10.68774818497546
$ws.Range("L17").Value = 11.50862537597219
$ws.Range("M17").Value = 17.50391347028441
$ws.Range("O17").Value = 14.58406179293042
$ws.Range("B18").Value = 19.26165950473131
$ws.Range("D18").Value = 3.71935206519174
$ws.Range("E18").Value = 35.38929918231891
$ws.Range("F18").Value = 16.68897657610724
$ws.Range("G18").Value = 16.12205991966125
$ws.Range("H18").Value = 10.70106625095675
$ws.Range("L18").Value = 11.42552756565282
$ws.Range("M18").Value = 17.42559760140519
$ws.Range("O18").Value = 14.60044797321583
$ws.Range("B19").Value = 19.22684300184318
$ws.Range("D19").Value = 3.71172751366831
$ws.Range("E19").Value = 35.30732384992424
$ws.Range("F19").Value = 16.68731428441866
$ws.Range("G19").Value = 16.11465863871257
$ws.Range("H19").Value = 10.70562762516655
$ws.Range("L19").Value = 11.39723683095509
$ws.Range("M19").Value = 17.39899639157301
$ws.Range("O19").Value = 14.60611128580807
$ws.Range("B20").Value = 19.38305364866742
$ws.Range("D20").Value = 3.745891906848273
$ws.Range("E20").Value = 35.67501890386495
$ws.Range("F20").Value = 16.69539505004424
$ws.Range("G20").Value = 16.14890221856367
$ws.Range("H20").Value = 10.68530808200319
$ws.Range("L20").Value = 11.52393074270045
$ws.Range("M20").Value = 17.51836723872271
$ws.Range("O20").Value = 14.58108391395721
$ws.Range("B21").Value = 19.89933092591266
$ws.Range("D21").Value = 3.85803697948716
$ws.Range("E21").Value = 36.88821592253397
$ws.Range("F21").Value = 16.73306038006442
$ws.Range("G21").Value = 16.28028202040133
$ws.Range("H21").Value = 10.62067518104426
$ws.Range("L21").Value = 11.93855469347232
$ws.Range("M21").Value = 17.91326525743871
$ws.Range("O21").Value = 14.50507119008138
$ws.Range("B22").Value = 20.22991202569708
$ws.Range("D22").Value = 3.929279332925172
$ws.Range("E22").Value = 37.66324253433623
$ws.Range("F22").Value = 16.76551775865496
$ws.Range("G22").Value = 16.37820182998934
$ws.Range("H22").Value = 10.58116097993368
$ws.Range("L22").Value = 12.20102606304448
$ws.Range("M22").Value = 18.1664454583445
$ws.Range("O22").Value = 14.46148041357025
$ws.Range("B23").Value = 20.05412788703011
$ws.Range("D23").Value = 3.891448572491954
$ws.Range("E23").Value = 37.25131123616232
$ws.Range("F23").Value = 16.74747742333644
$ws.Range("G23").Value = 16.32484572530806
$ws.Range("H23").Value = 10.60199990542066
$ws.Range("L23").Value = 12.06173520474785
$ws.Range("M23").Value = 18.0317867296089
$ws.Range("O23").Value = 14.48418422941188
$ws.Range("B24").Value = 19.37450635253461
$ws.Range("D24").Value = 3.744025484991351
$ws.Range("E24").Value = 35.6549068591841
$ws.Range("F24").Value = 16.69491176414521
$ws.Range("G24").Value = 16.14696007709908
$ws.Range("H24").Value = 10.68641028772914
$ws.Range("L24").Value = 11.51701414596675
$ws.Range("M24").Value = 17.51183436776092
$ws.Range("O24").Value = 14.58242809438329
$ws.Range("B25").Value = 18.61507401968716
$ws.Range("D25").Value = 3.57674860558277
$ws.Range("E25").Value = 33.8650957430238
$ws.Range("F25").Value = 16.67181246071895
$ws.Range("G25").Value = 16.00745527994198
$ws.Range("H25").Value = 10.78906802541786
$ws.Range("L25").Value = 10.89477466431073
$ws.Range("M25").Value = 16.93195607652495
$ws.Range("O25").Value = 14.71408605439039
